$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1678.9231
$ws.Range("I33").Value = 501.5
$ws.Range("K33").Value = 501.5
$ws.Range("M33").Value = -272.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1783.2142
$ws.Range("I43").Value = 1776.5555
$ws.Range("J43").Value = 1795.2
$ws.Range("K43").Value = 1776.5555
$ws.Range("L43").Value = 1795.2
$ws.Range("M43").Value = -1707.5555
$ws.Range("N43").Value = -1933.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 7932
$ws.Range("I76").Value = 6898
$ws.Range("J76").Value = 10000
$ws.Range("K76").Value = 6898
$ws.Range("L76").Value = 10000
$ws.Range("M76").Value = -6583
$ws.Range("N76").Value = -10630

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 7932
$ws.Range("I79").Value = 6898
$ws.Range("J79").Value = 10000
$ws.Range("K79").Value = 6898
$ws.Range("L79").Value = 10000
$ws.Range("M79").Value = -5806
$ws.Range("N79").Value = -12184

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 38967.5
$ws.Range("I92").Value = 1312.8667
$ws.Range("K92").Value = 1312.8667
$ws.Range("M92").Value = -64.86670000000004

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7643.206
$ws.Range("I116").Value = 6944.5264
$ws.Range("K116").Value = 6944.5264
$ws.Range("M116").Value = -3502.5264

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 34271.91
$ws.Range("I28").Value = 14998.5
$ws.Range("K28").Value = 14998.5
$ws.Range("M28").Value = -14806.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2646.6667
$ws.Range("I32").Value = 2617.0908
$ws.Range("K32").Value = 2617.0908
$ws.Range("M32").Value = -2330.0908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 38750
$ws.Range("I34").Value = 38000
$ws.Range("K34").Value = 38000
$ws.Range("M34").Value = -37729

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3980.6667
$ws.Range("J61").Value = 4539.9
$ws.Range("L61").Value = 4539.9
$ws.Range("N61").Value = -4963.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 34271.91
$ws.Range("I99").Value = 14998.5
$ws.Range("K99").Value = 14998.5
$ws.Range("M99").Value = -12003.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2940.125
$ws.Range("I110").Value = 2590.074
$ws.Range("K110").Value = 2590.074
$ws.Range("M110").Value = -545.0740000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1744.0834
$ws.Range("J122").Value = 1923.125
$ws.Range("L122").Value = 5769.375
$ws.Range("N122").Value = -10669.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3775.672
$ws.Range("I132").Value = 4341.125
$ws.Range("J132").Value = 2698.6191
$ws.Range("K132").Value = 13023.375
$ws.Range("L132").Value = 8095.8573
$ws.Range("M132").Value = -10493.375
$ws.Range("N132").Value = -13155.8573

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3980.6667
$ws.Range("J136").Value = 4539.9
$ws.Range("L136").Value = 13619.7
$ws.Range("N136").Value = -18719.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1167.4546
$ws.Range("I94").Value = 1159.9231
$ws.Range("K94").Value = 1159.9231
$ws.Range("M94").Value = -708.9231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3122.926
$ws.Range("I31").Value = 1352.8636
$ws.Range("K31").Value = 1352.8636
$ws.Range("M31").Value = -1057.8636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3122.926
$ws.Range("I34").Value = 1352.8636
$ws.Range("K34").Value = 1352.8636
$ws.Range("M34").Value = -1150.8636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2469
$ws.Range("I62").Value = 2589.25
$ws.Range("J62").Value = 1988
$ws.Range("K62").Value = 2589.25
$ws.Range("L62").Value = 1988
$ws.Range("M62").Value = -1965.25
$ws.Range("N62").Value = -3236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2469
$ws.Range("I65").Value = 2589.25
$ws.Range("J65").Value = 1988
$ws.Range("K65").Value = 12946.25
$ws.Range("L65").Value = 9940
$ws.Range("M65").Value = -9826.25
$ws.Range("N65").Value = -16180

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 193.72728
$ws.Range("I2").Value = 247.65218
$ws.Range("J2").Value = 69.7
$ws.Range("K2").Value = 1485.91308
$ws.Range("L2").Value = 418.2
$ws.Range("M2").Value = -1372.91308
$ws.Range("N2").Value = -644.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 42.77778
$ws.Range("J38").Value = 27.75
$ws.Range("L38").Value = 83.25
$ws.Range("N38").Value = -777.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1454.5938
$ws.Range("J68").Value = 1627.1666
$ws.Range("L68").Value = 4881.4998
$ws.Range("N68").Value = -6503.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1454.5938
$ws.Range("J71").Value = 1627.1666
$ws.Range("L71").Value = 14644.4994
$ws.Range("N71").Value = -22756.4994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1420.1666
$ws.Range("I92").Value = 658.8333
$ws.Range("K92").Value = 1976.4999
$ws.Range("M92").Value = -728.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 15494.75
$ws.Range("I116").Value = 15494.75
$ws.Range("K116").Value = 46484.25
$ws.Range("M116").Value = -43042.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 941.6667
$ws.Range("I119").Value = 941.6667
$ws.Range("K119").Value = 2825.0001
$ws.Range("M119").Value = 2012.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3050.8572
$ws.Range("I80").Value = 3104.182
$ws.Range("J80").Value = 2855.3333
$ws.Range("K80").Value = 3104.182
$ws.Range("L80").Value = 2855.3333
$ws.Range("M80").Value = -2106.182
$ws.Range("N80").Value = -4851.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3050.8572
$ws.Range("I83").Value = 3104.182
$ws.Range("J83").Value = 2855.3333
$ws.Range("K83").Value = 15520.91
$ws.Range("L83").Value = 14276.6665
$ws.Range("M83").Value = -10528.91
$ws.Range("N83").Value = -24260.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8125.3335
$ws.Range("I102").Value = 9524.666999999999
$ws.Range("J102").Value = 6026.3335
$ws.Range("K102").Value = 9524.666999999999
$ws.Range("L102").Value = 6026.3335
$ws.Range("M102").Value = -7902.666999999999
$ws.Range("N102").Value = -9270.333500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 566.2778
$ws.Range("I107").Value = 515.1818
$ws.Range("J107").Value = 646.5714
$ws.Range("K107").Value = 515.1818
$ws.Range("L107").Value = 646.5714
$ws.Range("M107").Value = 1404.8182
$ws.Range("N107").Value = -4486.5714

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 24853.564
$ws.Range("I122").Value = 35142.355
$ws.Range("J122").Value = 3590.0667
$ws.Range("K122").Value = 105427.065
$ws.Range("L122").Value = 10770.2001
$ws.Range("M122").Value = -102977.065
$ws.Range("N122").Value = -15670.2001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 49999.5
$ws.Range("I25").Value = 49999
$ws.Range("K25").Value = 49999
$ws.Range("M25").Value = -49769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I100").Value = 71429224
$ws.Range("K100").Value = 142858448
$ws.Range("M100").Value = -142857907

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1787.9584
$ws.Range("I113").Value = 1017.86957
$ws.Range("K113").Value = 3053.60871
$ws.Range("M113").Value = -883.60871

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4198.923
$ws.Range("I132").Value = 4365.3335
$ws.Range("K132").Value = 13096.0005
$ws.Range("M132").Value = -10566.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3078573.8
$ws.Range("I136").Value = 4526144.5
$ws.Range("K136").Value = 13578433.5
$ws.Range("M136").Value = -13575883.5
